$d = $word.ActiveDocument

# The "Bibliografia" paragraph is currently one long run of text with all
# references run together. Split it into one run per reference by
# inserting a manual line break (the "^l" wildcard code, i.e. <w:br/>)
# right after each reference ends and before the next one begins.
# Each Find/Replace below is scoped to a unique substring so it only
# touches the intended boundary.

$d.Content.Find.Execute(
    "2005. BENNETT", $true, $false, $false, $false, $false,
    $true, 1, $false, "2005. ^lBENNETT", 2) | Out-Null

$d.Content.Find.Execute(
    "McGraw-Hill. KREITH", $true, $false, $false, $false, $false,
    $true, 1, $false, "McGraw-Hill. ^lKREITH", 2) | Out-Null

$d.Content.Find.Execute(
    "2003.BIRD", $true, $false, $false, $false, $false,
    $true, 1, $false, "2003.^lBIRD", 2) | Out-Null

$d.Content.Find.Execute(
    "2004. FOX", $true, $false, $false, $false, $false,
    $true, 1, $false, "2004. ^lFOX", 2) | Out-Null

$d.Content.Find.Execute(
    "2001. SISSOM", $true, $false, $false, $false, $false,
    $true, 1, $false, "2001. ^lSISSOM", 2) | Out-Null

$d.Content.Find.Execute(
    "1988. HOLMAN", $true, $false, $false, $false, $false,
    $true, 1, $false, "1988. ^lHOLMAN", 2) | Out-Null

$d.Content.Find.Execute(
    "1983. POIRIER", $true, $false, $false, $false, $false,
    $true, 1, $false, "1983. ^lPOIRIER", 2) | Out-Null

$d.Content.Find.Execute(
    "1994.GASKELL", $true, $false, $false, $false, $false,
    $true, 1, $false, "1994.^lGASKELL", 2) | Out-Null

$d.Content.Find.Execute(
    "1991. SZEKELY", $true, $false, $false, $false, $false,
    $true, 1, $false, "1991. ^lSZEKELY", 2) | Out-Null
